$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New UoM rows appended to the table (Box, Roll, Set)
$data = @(
    @("Box",  "Buah", "Bigger than the reference Unit of Measure",   1, "Yes"),
    @("Roll", "Buah", "Bigger than the reference Unit of Measure",   1, "Yes"),
    @("Set",  "Buah", "Reference Unit of Measure for this category", 1, "Yes")
)

$startRow = 17
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]

    # Match the wrap/vertical-center formatting used by the rest of the table body
    $rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 5))
    $rowRange.VerticalAlignment = -4108
    $rowRange.WrapText = $true
}

# Grow the table to cover the new rows
[void]$ws.ListObjects.Item("Table").Resize($ws.Range("A1:E19"))

# Selection ends up on the last new cell, matching the saved view
[void]$ws.Range("E19").Select()
